$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 698.47
$ws.Range("I15").Value = 698.47
$ws.Range("K15").Value = 2095.41
$ws.Range("M15").Value = -1926.41
$ws.Range("H17").Value = 1277.1111
$ws.Range("J17").Value = 1299.8823
$ws.Range("L17").Value = 3899.6469
$ws.Range("N17").Value = -4235.6469
$ws.Range("H58").Value = 2545.6667
$ws.Range("I58").Value = 789.375
$ws.Range("K58").Value = 2368.125
$ws.Range("M58").Value = -2218.125
$ws.Range("H112").Value = 3162.0688
$ws.Range("J112").Value = 3404
$ws.Range("L112").Value = 10212
$ws.Range("N112").Value = -12428
$ws.Range("H132").Value = 6997.5415
$ws.Range("I132").Value = 5542.15
$ws.Range("J132").Value = 14274.5
$ws.Range("K132").Value = 16626.45
$ws.Range("L132").Value = 42823.5
$ws.Range("M132").Value = -14096.45
$ws.Range("N132").Value = -47883.5
$ws.Range("H137").Value = 2359.44
$ws.Range("I137").Value = 3331.2307
$ws.Range("K137").Value = 9993.6921
$ws.Range("M137").Value = -7443.6921
$ws.Range("H138").Value = 3628.9019
$ws.Range("I138").Value = 2483.52
$ws.Range("J138").Value = 4730.231
$ws.Range("K138").Value = 7450.559999999999
$ws.Range("L138").Value = 14190.693
$ws.Range("M138").Value = -2310.559999999999
$ws.Range("N138").Value = -24470.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3481.25
$ws.Range("I61").Value = 3562.2856
$ws.Range("J61").Value = 2914
$ws.Range("K61").Value = 3562.2856
$ws.Range("L61").Value = 2914
$ws.Range("M61").Value = -3350.2856
$ws.Range("N61").Value = -3338
$ws.Range("H76").Value = 18500
$ws.Range("J76").Value = 18500
$ws.Range("L76").Value = 18500
$ws.Range("N76").Value = -19176
$ws.Range("H79").Value = 18500
$ws.Range("J79").Value = 18500
$ws.Range("L79").Value = 18500
$ws.Range("N79").Value = -20840
$ws.Range("H132").Value = 6660.1
$ws.Range("I132").Value = 6172.759
$ws.Range("J132").Value = 7333.095
$ws.Range("K132").Value = 18518.277
$ws.Range("L132").Value = 21999.285
$ws.Range("M132").Value = -15988.277
$ws.Range("N132").Value = -27059.285
$ws.Range("H136").Value = 3481.25
$ws.Range("I136").Value = 3562.2856
$ws.Range("J136").Value = 2914
$ws.Range("K136").Value = 10686.8568
$ws.Range("L136").Value = 8742
$ws.Range("M136").Value = -8136.856800000001
$ws.Range("N136").Value = -13842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2599.7886
$ws.Range("I134").Value = 1739.6666
$ws.Range("J134").Value = 3528.72
$ws.Range("K134").Value = 5218.9998
$ws.Range("L134").Value = 10586.16
$ws.Range("M134").Value = -2683.9998
$ws.Range("N134").Value = -15656.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2954
$ws.Range("I31").Value = 1793.5227
$ws.Range("J31").Value = 5641.421
$ws.Range("K31").Value = 1793.5227
$ws.Range("L31").Value = 5641.421
$ws.Range("M31").Value = -1498.5227
$ws.Range("N31").Value = -6231.421
$ws.Range("H34").Value = 2954
$ws.Range("I34").Value = 1793.5227
$ws.Range("J34").Value = 5641.421
$ws.Range("K34").Value = 1793.5227
$ws.Range("L34").Value = 5641.421
$ws.Range("M34").Value = -1591.5227
$ws.Range("N34").Value = -6045.421
$ws.Range("H74").Value = 33300
$ws.Range("J74").Value = 33300
$ws.Range("L74").Value = 33300
$ws.Range("N74").Value = -35048
$ws.Range("H77").Value = 33300
$ws.Range("J77").Value = 33300
$ws.Range("L77").Value = 99900
$ws.Range("N77").Value = -108636
$ws.Range("H132").Value = 1591
$ws.Range("I132").Value = 1034.6857
$ws.Range("J132").Value = 2262.4138
$ws.Range("K132").Value = 3104.0571
$ws.Range("L132").Value = 6787.241399999999
$ws.Range("M132").Value = -574.0571
$ws.Range("N132").Value = -11847.2414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2261.28
$ws.Range("I22").Value = 1795.7778
$ws.Range("J22").Value = 2523.125
$ws.Range("K22").Value = 5387.3334
$ws.Range("L22").Value = 7569.375
$ws.Range("M22").Value = -5218.3334
$ws.Range("N22").Value = -7907.375
$ws.Range("H27").Value = 2261.28
$ws.Range("I27").Value = 1795.7778
$ws.Range("J27").Value = 2523.125
$ws.Range("K27").Value = 5387.3334
$ws.Range("L27").Value = 7569.375
$ws.Range("M27").Value = -5285.3334
$ws.Range("N27").Value = -7773.375
$ws.Range("H34").Value = 12067.444
$ws.Range("J34").Value = 21501.2
$ws.Range("L34").Value = 64503.60000000001
$ws.Range("N34").Value = -64671.60000000001
$ws.Range("H39").Value = 1012.7143
$ws.Range("J39").Value = 1103.1111
$ws.Range("L39").Value = 3309.3333
$ws.Range("N39").Value = -3897.3333
$ws.Range("H55").Value = 14374.934
$ws.Range("J55").Value = 15366
$ws.Range("L55").Value = 46098
$ws.Range("N55").Value = -46452
$ws.Range("H59").Value = 2700.25
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2700.25
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 8100.75
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -9180.75
$ws.Range("H68").Value = 1011.5
$ws.Range("I68").Value = 870.2857
$ws.Range("K68").Value = 2610.8571
$ws.Range("M68").Value = -1799.8571
$ws.Range("H71").Value = 1011.5
$ws.Range("I71").Value = 870.2857
$ws.Range("K71").Value = 7832.571300000001
$ws.Range("M71").Value = -3776.571300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7309.3335
$ws.Range("I122").Value = 10202.071
$ws.Range("J122").Value = 3259.5
$ws.Range("K122").Value = 30606.213
$ws.Range("L122").Value = 9778.5
$ws.Range("M122").Value = -28156.213
$ws.Range("N122").Value = -14678.5
$ws.Range("H132").Value = 2397.182
$ws.Range("I132").Value = 1924.6086
$ws.Range("J132").Value = 2736.8438
$ws.Range("K132").Value = 5773.825800000001
$ws.Range("L132").Value = 8210.5314
$ws.Range("M132").Value = -3243.825800000001
$ws.Range("N132").Value = -13270.5314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17357.463
$ws.Range("I132").Value = 22216.418
$ws.Range("J132").Value = 7459.593
$ws.Range("K132").Value = 66649.254
$ws.Range("L132").Value = 22378.779
$ws.Range("M132").Value = -64119.254
$ws.Range("N132").Value = -27438.779
$ws.Range("H136").Value = 1660.2174
$ws.Range("I136").Value = 1396.6177
$ws.Range("J136").Value = 2407.0833
$ws.Range("K136").Value = 4189.8531
$ws.Range("L136").Value = 7221.249899999999
$ws.Range("M136").Value = -1639.8531
$ws.Range("N136").Value = -12321.2499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1618.2
$ws.Range("I122").Value = 1750.2632
$ws.Range("K122").Value = 5250.7896
$ws.Range("M122").Value = -2800.7896
